{"js": "// Update the date line and the 25 multiplication problems in the table\n// to the new values from the latest generated output.\nconst replacements = [\n  [\"2025-12-04 Thursday\", \"2025-12-05 Friday\"],\n  [\"881\u00d76=\", \"324\u00d76=\"],\n  [\"848\u00d78=\", \"866\u00d77=\"],\n  [\"220\u00d75=\", \"877\u00d76=\"],\n  [\"302\u00d73=\", \"122\u00d77=\"],\n  [\"284\u00d79=\", \"735\u00d78=\"],\n  [\"665\u00d76=\", \"254\u00d79=\"],\n  [\"500\u00d72=\", \"976\u00d76=\"],\n  [\"918\u00d76=\", \"566\u00d72=\"],\n  [\"831\u00d76=\", \"672\u00d72=\"],\n  [\"803\u00d72=\", \"229\u00d74=\"],\n  [\"271\u00d77=\", \"982\u00d77=\"],\n  [\"345\u00d75=\", \"121\u00d76=\"],\n  [\"431\u00d74=\", \"373\u00d78=\"],\n  [\"869\u00d77=\", \"163\u00d77=\"],\n  [\"411\u00d72=\", \"359\u00d77=\"],\n  [\"479\u00d76=\", \"667\u00d72=\"],\n  [\"999\u00d79=\", \"750\u00d74=\"],\n  [\"738\u00d75=\", \"496\u00d73=\"],\n  [\"482\u00d79=\", \"432\u00d76=\"],\n  [\"196\u00d73=\", \"777\u00d76=\"],\n  [\"155\u00d75=\", \"347\u00d77=\"],\n  [\"488\u00d72=\", \"280\u00d72=\"],\n  [\"821\u00d79=\", \"455\u00d76=\"],\n  [\"923\u00d77=\", \"649\u00d79=\"],\n  [\"763\u00d74=\", \"691\u00d72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 multiplication problems in the table\n# to the new values from the latest generated output.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-12-04 Thursday\", \"2025-12-05 Friday\"),\n    @(\"881\u00d76=\", \"324\u00d76=\"),\n    @(\"848\u00d78=\", \"866\u00d77=\"),\n    @(\"220\u00d75=\", \"877\u00d76=\"),\n    @(\"302\u00d73=\", \"122\u00d77=\"),\n    @(\"284\u00d79=\", \"735\u00d78=\"),\n    @(\"665\u00d76=\", \"254\u00d79=\"),\n    @(\"500\u00d72=\", \"976\u00d76=\"),\n    @(\"918\u00d76=\", \"566\u00d72=\"),\n    @(\"831\u00d76=\", \"672\u00d72=\"),\n    @(\"803\u00d72=\", \"229\u00d74=\"),\n    @(\"271\u00d77=\", \"982\u00d77=\"),\n    @(\"345\u00d75=\", \"121\u00d76=\"),\n    @(\"431\u00d74=\", \"373\u00d78=\"),\n    @(\"869\u00d77=\", \"163\u00d77=\"),\n    @(\"411\u00d72=\", \"359\u00d77=\"),\n    @(\"479\u00d76=\", \"667\u00d72=\"),\n    @(\"999\u00d79=\", \"750\u00d74=\"),\n    @(\"738\u00d75=\", \"496\u00d73=\"),\n    @(\"482\u00d79=\", \"432\u00d76=\"),\n    @(\"196\u00d73=\", \"777\u00d76=\"),\n    @(\"155\u00d75=\", \"347\u00d77=\"),\n    @(\"488\u00d72=\", \"280\u00d72=\"),\n    @(\"821\u00d79=\", \"455\u00d76=\"),\n    @(\"923\u00d77=\", \"649\u00d79=\"),\n    @(\"763\u00d74=\", \"691\u00d72=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $null, $null, $find.Forward, $find.Wrap, $find.Format, $find.Replacement.Text, 2)  # wdReplaceAll\n}\n"}
